$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "add more tag res" - append a third quiz-answer row (tag = "ถูกสาม") whose
# response is the congratulations message (with the prize link) shown once
# the player has answered all three trivia questions correctly.
$ws.Range("A45").Value = "ถูกสาม"
$ws.Range("B45").Value = "ยินดีด้วยค้าบเตง ตอบถูกหมดเยย รับไปสำหรับรางวัลของคนเก่ง <3`nhttps://www.youtube.com/watch?v=dQw4w9WgXcQ"

# Match the formatting used by the other multi-line "response" cells: wrapped
# text with a 30pt row height.
$ws.Range("B45").WrapText = $true
$ws.Rows.Item(45).RowHeight = 30

# Scroll the view down a bit (as happened after the row was appended) and
# leave the same cell selected that the source workbook ended up with.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B48").Select()
